$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update assistant message texts (shared strings used in column B)
$ws.Range("B3").Value = "I am Cuong, your AI assistant. How can I help you today?"
$ws.Range("B4").Value = "Hello Cuong! I'm Minh, your AI assistant. How can I assist you today?"
$ws.Range("B5").Value = "Hello Minh! It seems like there's a little mix-up, but I'm here to assist you. How can I help you today?"
$ws.Range("B6").Value = "Thank you, Cuong! I appreciate your willingness to assist. If you have any questions or need information on a specific topic, feel free to ask!"

$ws.Range("B10").Value = "Hi Cuong! I'm here to assist you as well. How can I help you today?"
$ws.Range("B11").Value = "Hi there! I appreciate your offer to help. I'm here to provide information and answer any questions you might have. What would you like to know or discuss today?"
$ws.Range("B12").Value = "Thank you, Cuong! I appreciate your willingness to help. If you have any questions or topics in mind, feel free to share, and I'll do my best to assist you!"

# Update Response_Time values
$ws.Range("C2").Value = 1.115615606307983
$ws.Range("C3").Value = 0.7531790733337402
$ws.Range("C4").Value = 0.9417111873626709
$ws.Range("C5").Value = 1.120330810546875

$ws.Range("C8").Value = 0.8064866065979004
$ws.Range("C9").Value = 0.8226840496063232
$ws.Range("C10").Value = 1.499311447143555
$ws.Range("C11").Value = 1.115266799926758
